$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Build the locale suffixes used in the "name" column, e.g. " – Español" / " – Inglés"
$dash = [char]0x2013
$suffixEs = " " + $dash + " Espa" + [char]0xF1 + "ol"
$suffixEn = " " + $dash + " Ingl" + [char]0xE9 + "s"

# --- Update Hoja2 (data table) ---
# Column B (name): strip the trailing " – Español"/" – Inglés" suffix.
# Column C (set): replace the old "SV ..." product family text with the category code (column E).
for ($r = 2; $r -le 10; $r++) {
    $name = $ws2.Cells.Item($r, 2).Value2
    if ($name.EndsWith($suffixEs)) {
        $ws2.Cells.Item($r, 2).Value = $name.Substring(0, $name.Length - $suffixEs.Length)
    } elseif ($name.EndsWith($suffixEn)) {
        $ws2.Cells.Item($r, 2).Value = $name.Substring(0, $name.Length - $suffixEn.Length)
    }

    $category = $ws2.Cells.Item($r, 5).Value2
    $ws2.Cells.Item($r, 3).Value = $category
}

# Column B is narrower now that the suffix was removed; resize it to fit the new content.
$ws2.Columns.Item(2).ColumnWidth = 52.17

# --- Update Hoja1 (plain pasted-values copy of Hoja2 column O) ---
# Refresh rows 2..9 with the (now recalculated) concatenated text from Hoja2,
# then drop the old row 10 (etb-surging-sparks-en), which is no longer listed.
for ($r = 2; $r -le 9; $r++) {
    $ws1.Cells.Item($r, 1).Value = $ws2.Cells.Item($r, 15).Value2
}

$ws1.Rows.Item(10).Delete()

$ws1.Range("A1:A9").Select() | Out-Null
$ws1.Activate() | Out-Null
